$d = $word.ActiveDocument

# --- 1. Heading: Review number / title ---
$d.Content.Find.Execute(
    "Review 180: [Short] Learning From Mistakes Makes LLM Better Reasoner",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 179: [Short] UNLEASHING THE POWER OF PRE-TRAINED LANGUAGE MODELS FOR OFFLINE REINFORCEMENT LEARNING",
    2)

# --- 2. Bold "Paper:" link line ---
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2310.20689v4",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2310.20587v5",
    2)

# --- 3. huggingface link line ---
$d.Content.Find.Execute(
    "https://huggingface.co/papers/2310.20689",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://huggingface.co/papers/2310.20587",
    2)

# --- 4. Insert two new empty "Normal" paragraphs plus a new Hebrew paragraph
#        right after the huggingface link paragraph (old paragraph 4). ---
$hfPara = $d.Paragraphs.Item(4)
$hfPara.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item(5)
$blank1.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item(6)
$newPara = $blank2.Range.InsertParagraphAfter()
$introPara = $d.Paragraphs.Item(7)
$introPara.Range.InsertAfter("המאמר הזה משך את עיני כי אני מאוד אוהב שילובים של מודלי שפה (וגם מודלי דיפוזיה) למשימות מהעולם של למידה עם חיזוקים (Reinforcement Learning). המאמר הזה עשה את זה בצורה מאוד אלגנטית כל כך אחרי שרפרפתי בו 5 דקות לא הבנתי חשבתי שזה די בלתי אפשרי. אבל אחרי הקצת צללתי לעומק הבנתי שכדאי לתת לזה צ'אנס ולסקור אותו במדורנו.")

# --- 5. Replace the long two-part (with line breaks) paragraph with the new single paragraph ---
# This was old paragraph 6 ("טוב, שוב בשבת..."), now shifted to index 9.
# (Clear the paragraph's contents first, then insert fresh text, so the
#  run doesn't inherit the old xml:space="preserve" run formatting.)
$bodyPara1 = $d.Paragraphs.Item(9)
$r1 = $bodyPara1.Range
$r1.End = $r1.End - 1
$r1.Text = ""
$r1.InsertAfter("אז מה המאמר עשה בעצם? בגדול הם לקחו מודל שפה וטייבו (finetuned) אותו לבצע למשימות של RL. כלומר בהינתן של פעולות ומצבים קודמים המטרה של המודל היא לחזות את הפעולה הבאה. במקרה הזה מדובר באופליין RL כלומר המטרה של המודל היא לחקות כמה שיותר טוב את הפעולות המוצלחות מהדאטהסט (בהינתן הפעולות והמצבים הקודמים). במשימות שנדונו במאמר הפעולות מתוארות בצורה מילולית.")

# --- 6. Remove the now-empty paragraph that followed it (old paragraph 7, now index 10) ---
$d.Paragraphs.Item(10).Range.Delete()

# --- 7. Replace the final body paragraph's text (old paragraph 8, now index 10) ---
$bodyPara2 = $d.Paragraphs.Item(10)
$bodyPara2.Range.Text = "כבר מריחים את מודלי השפה מתקרבים? אוקיי, קודם המחברים לקחו מודל שפה מאומן (GPT2) וטייבו אותו על הדאטהסט הנקרא WiKiText. בשלב השני מוסיפים למודל שפה כמה שכבות של MLP ומאמנים אותו על הדאטה של המשימה (נגיד משחק אטארי) ובנוסף מכיילים מודל שפה עם LoRA (זוכרים מה זה?). תוך כדי התהליך הזה מוסיפים איבר רגולריזציה המכיל לוס על הדאטה של WikiText כנראה כדי לגרום למודל לא לשכוח את המיומנויות הקודמות שלו. וזה וזה מקבלים מודל ל-RL כלומר decision transformer עם ביצועים טובים."

# --- 8. Remove the trailing empty paragraph (old paragraph 9, now the last one) ---
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
